$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.942.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.914.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.398.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.923.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.917.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "430.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.679"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +5.50%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.704.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("E51").Value = "  +1.45%  "
